# Arjun Committed Code on 10-30-2019
# Populate additional BAN / subscriber test-data rows on INPUT_SHEET.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INPUT_SHEET")

# Columns A (BAN) and B (SM_USER) already carry a Text ("@") column style,
# so plain .Value assignment is enough to keep numeric-looking ids (with
# leading zeros) stored as text rather than being coerced to numbers.
$rows = @(
    @{ Row = 2;  BAN = "107198053"; User = "ONREG-20535"; Id = "08072919121"; Status = "ACTIVE" },
    @{ Row = 3;  BAN = "124473304"; User = "ONREG-18520"; Id = "08317211011"; Status = "ACTIVE" },
    @{ Row = 4;  BAN = "169068211"; User = "ONREG-19643"; Id = "03085439121"; Status = "ACTIVE" },
    @{ Row = 5;  BAN = "225356050"; User = "ONREG-22679"; Id = "21545071011"; Status = "ACTIVE" },
    @{ Row = 6;  BAN = "244834015"; User = "ONREG-22710"; Id = "00130829021"; Status = "ACTIVE" },
    @{ Row = 7;  BAN = "114478789"; User = "ONREG-25767"; Id = "13430032021"; Status = "ACTIVE" },
    @{ Row = 8;  BAN = "172500199"; User = "ONREG-21172"; Id = "35296929121"; Status = "ACTIVE" },
    @{ Row = 9;  BAN = "205653581"; User = "ONREG-26222"; Id = "15456191021"; Status = "ACTIVE" },
    @{ Row = 10; BAN = "";          User = "";            Id = "NA";          Status = "NA" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.BAN
    $ws.Cells.Item($r.Row, 2).Value = $r.User
    $ws.Cells.Item($r.Row, 3).Value = $r.Id
    $ws.Cells.Item($r.Row, 4).Value = $r.Status
}

# Column B (SM_USER) widened to fit the longer ONREG ids.
$ws.Columns.Item(2).ColumnWidth = 15.42578125

# Selection lands on the last data row, matching the author's final click.
$ws.Range("A9:XFD9").Select()
